$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns O (Operacion) and P (Zona)
$ws.Range("O1").Value = "Operacion"
$ws.Range("P1").Value = "Zona"

# Copy header style from an existing header cell (N1) onto the two new header cells
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)

# Populate data rows 2-36 for columns O (Operacion) and P (Zona)
$ws.Range("O2").Value = "San Telmo"
$ws.Range("P2").Value = "Capital Sur"
$ws.Range("O3").Value = "Recoleta"
$ws.Range("P3").Value = "Capital Sur"
$ws.Range("O4").Value = "Saavedra"
$ws.Range("P4").Value = "Capital Norte"
$ws.Range("O5").Value = "Colegiales"
$ws.Range("P5").Value = "Capital Norte"
$ws.Range("O6").Value = "San Telmo"
$ws.Range("P6").Value = "Capital Sur"
$ws.Range("O7").Value = "Recoleta"
$ws.Range("P7").Value = "Capital Sur"
$ws.Range("O8").Value = "Paternal"
$ws.Range("P8").Value = "Capital Norte"
$ws.Range("O9").Value = "San Telmo"
$ws.Range("P9").Value = "Capital Sur"
$ws.Range("O10").Value = "San Telmo"
$ws.Range("P10").Value = "Capital Sur"
$ws.Range("O11").Value = "Almagro"
$ws.Range("P11").Value = "Capital Sur"
$ws.Range("O12").Value = "San Telmo"
$ws.Range("P12").Value = "Capital Sur"
$ws.Range("O13").Value = "San Telmo"
$ws.Range("P13").Value = "Capital Sur"
$ws.Range("O14").Value = "Colegiales"
$ws.Range("P14").Value = "Capital Norte"
$ws.Range("O15").Value = "Recoleta"
$ws.Range("P15").Value = "Capital Sur"
$ws.Range("O16").Value = "Recoleta"
$ws.Range("P16").Value = "Capital Sur"
$ws.Range("O17").Value = "Devoto"
$ws.Range("P17").Value = "Capital Norte"
$ws.Range("O18").Value = "Palermo"
$ws.Range("P18").Value = "Capital Sur"
$ws.Range("O19").Value = "Saavedra"
$ws.Range("P19").Value = "Capital Norte"
$ws.Range("O20").Value = "Saavedra"
$ws.Range("P20").Value = "Capital Norte"
$ws.Range("O21").Value = "Paternal"
$ws.Range("P21").Value = "Capital Norte"
$ws.Range("O22").Value = "Saavedra"
$ws.Range("P22").Value = "Capital Norte"
$ws.Range("O23").Value = "Saavedra"
$ws.Range("P23").Value = "Capital Norte"
$ws.Range("O24").Value = "Saavedra"
$ws.Range("P24").Value = "Capital Norte"
$ws.Range("O25").Value = "Saavedra"
$ws.Range("P25").Value = "Capital Norte"
$ws.Range("O26").Value = "Paternal"
$ws.Range("P26").Value = "Capital Norte"
$ws.Range("O27").Value = "San Telmo"
$ws.Range("P27").Value = "Capital Sur"
$ws.Range("O28").Value = "Paternal"
$ws.Range("P28").Value = "Capital Norte"
$ws.Range("O29").Value = "Almagro"
$ws.Range("P29").Value = "Capital Sur"
$ws.Range("O30").Value = "Devoto"
$ws.Range("P30").Value = "Capital Norte"
$ws.Range("O31").Value = "Devoto"
$ws.Range("P31").Value = "Capital Norte"
$ws.Range("O32").Value = "Paternal"
$ws.Range("P32").Value = "Capital Norte"
$ws.Range("O33").Value = "Colegiales"
$ws.Range("P33").Value = "Capital Norte"
$ws.Range("O34").Value = "Paternal"
$ws.Range("P34").Value = "Capital Norte"
$ws.Range("O35").Value = "Almagro"
$ws.Range("P35").Value = "Capital Sur"
$ws.Range("O36").Value = "San Telmo"
$ws.Range("P36").Value = "Capital Sur"
